$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date range in header (I1, K1): 14-10-2025 -> 15-10-2025
$ws.Range("I1").Value = "15-10-2025 00:00:00"
$ws.Range("K1").Value = "15-10-2025 00:00:00"

# Rows 161,162,163: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B161").Value = 53925
$ws.Range("C161").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D161").Value = 66.44
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44
$ws.Range("B162").Value = 64350
$ws.Range("C162").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D162").Value = 66.44
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 101
$ws.Range("G162").Value = 6710.44
$ws.Range("B163").Value = 57756
$ws.Range("C163").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D163").Value = 66.44
$ws.Range("E163").Value = 79.37
$ws.Range("F163").Value = -100
$ws.Range("G163").Value = -6644

# Rows 279,280: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B279").Value = 64973
$ws.Range("C279").Value = "HIM-GENTLE BABY SOAP 75G"
$ws.Range("D279").Value = 33.3
$ws.Range("E279").Value = 35.4
$ws.Range("F279").Value = 150
$ws.Range("G279").Value = 4995
$ws.Range("B280").Value = 48706
$ws.Range("C280").Value = "HIM-GENTLE BABY SOAP 75G"
$ws.Range("D280").Value = 33.3
$ws.Range("E280").Value = 39.8
$ws.Range("F280").Value = -144
$ws.Range("G280").Value = -4795.2

# Rows 313,314: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B313").Value = 57854
$ws.Range("C313").Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Range("D313").Value = 305.84
$ws.Range("E313").Value = 325.16
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.6799999999999
$ws.Range("B314").Value = 62997
$ws.Range("C314").Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Range("D314").Value = 305.84
$ws.Range("E314").Value = 325.16
$ws.Range("F314").Value = 72
$ws.Range("G314").Value = 22020.48

# Rows 316,317: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B316").Value = 63565
$ws.Range("C316").Value = "HUL-Bru Inst Poly 50g"
$ws.Range("D316").Value = 102.71
$ws.Range("E316").Value = 109.19
$ws.Range("F316").Value = 60
$ws.Range("G316").Value = 6162.6
$ws.Range("B317").Value = 57077
$ws.Range("C317").Value = "HUL-Bru Inst Poly 50g"
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08

# Rows 351,352: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B351").Value = 63531
$ws.Range("C351").Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Range("D351").Value = 143.48
$ws.Range("E351").Value = 152.53
$ws.Range("F351").Value = 80
$ws.Range("G351").Value = 11478.4
$ws.Range("B352").Value = 63571
$ws.Range("C352").Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Range("D352").Value = 143.48
$ws.Range("E352").Value = 152.53
$ws.Range("F352").Value = 27
$ws.Range("G352").Value = 3873.96

# Rows 355,356: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B355").Value = 63510
$ws.Range("C355").Value = "HUL-knorr schezwan 200g pch"
$ws.Range("D355").Value = 47.64
$ws.Range("E355").Value = 50.66
$ws.Range("F355").Value = 167
$ws.Range("G355").Value = 7955.88
$ws.Range("B356").Value = 55356
$ws.Range("C356").Value = "HUL-knorr schezwan 200g pch"
$ws.Range("D356").Value = 47.64
$ws.Range("E356").Value = 54.04
$ws.Range("F356").Value = -158
$ws.Range("G356").Value = -7527.12

# Rows 375,376: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B375").Value = 61605
$ws.Range("C375").Value = "HUL-lux advanced eventoned glow 4x100"
$ws.Range("D375").Value = 111.96
$ws.Range("E375").Value = 133.78
$ws.Range("F375").Value = -13
$ws.Range("G375").Value = -1455.48
$ws.Range("B376").Value = 63563
$ws.Range("C376").Value = "HUL-lux advanced eventoned glow 4x100"
$ws.Range("D376").Value = 111.96
$ws.Range("E376").Value = 119.04
$ws.Range("F376").Value = 15
$ws.Range("G376").Value = 1679.4

# Rows 382,383: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B382").Value = 63560
$ws.Range("C382").Value = "Hul-pears pure and gentle 3x125 gm"
$ws.Range("D382").Value = 126.86
$ws.Range("E382").Value = 134.87
$ws.Range("F382").Value = 104
$ws.Range("G382").Value = 13193.44
$ws.Range("B383").Value = 60325
$ws.Range("C383").Value = "Hul-pears pure and gentle 3x125 gm"
$ws.Range("D383").Value = 126.86
$ws.Range("E383").Value = 151.57
$ws.Range("F383").Value = -102
$ws.Range("G383").Value = -12939.72

# Rows 419,420: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B419").Value = 57856
$ws.Range("C419").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("D419").Value = 171.33
$ws.Range("E419").Value = 204.69
$ws.Range("F419").Value = 2
$ws.Range("G419").Value = 342.66
$ws.Range("B420").Value = 63007
$ws.Range("C420").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("D420").Value = 171.33
$ws.Range("E420").Value = 204.69
$ws.Range("F420").Value = 984
$ws.Range("G420").Value = 168588.72

# Rows 431,432: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("D431").Value = 59.47
$ws.Range("E431").Value = 71.05
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("D432").Value = 59.47
$ws.Range("E432").Value = 71.05
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47

# Rows 536,537: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B536").Value = 47097
$ws.Range("C536").Value = "KUS-Floor Wiper"
$ws.Range("D536").Value = 112.28
$ws.Range("E536").Value = 134.16
$ws.Range("F536").Value = 15
$ws.Range("G536").Value = 1684.2
$ws.Range("B537").Value = 58047
$ws.Range("C537").Value = "KUS-Floor Wiper"
$ws.Range("D537").Value = 105.54
$ws.Range("E537").Value = 126.1
$ws.Range("F537").Value = 54
$ws.Range("G537").Value = 5699.16

# Rows 579,580: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B579").Value = 53757
$ws.Range("C579").Value = "CRE-Bourbon 100gm"
$ws.Range("D579").Value = 13.45
$ws.Range("E579").Value = 16.08
$ws.Range("F579").Value = -159
$ws.Range("G579").Value = -2138.55
$ws.Range("B580").Value = 65069
$ws.Range("C580").Value = "CRE-Bourbon 100gm"
$ws.Range("D580").Value = 13.45
$ws.Range("E580").Value = 14.3
$ws.Range("F580").Value = 172
$ws.Range("G580").Value = 2313.4

# Rows 583,584: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B583").Value = 53263
$ws.Range("C583").Value = "CRE-Butter cremfills 100gm"
$ws.Range("D583").Value = 12.81
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29
$ws.Range("B584").Value = 65066
$ws.Range("C584").Value = "CRE-Butter cremfills 100gm"
$ws.Range("D584").Value = 12.81
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 313
$ws.Range("G584").Value = 4009.53

# Rows 586,587: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B586").Value = 64915
$ws.Range("C586").Value = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Range("D586").Value = 19.73
$ws.Range("E586").Value = 20.98
$ws.Range("F586").Value = 40
$ws.Range("G586").Value = 789.2
$ws.Range("B587").Value = 45695
$ws.Range("C587").Value = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Range("D587").Value = 19.73
$ws.Range("E587").Value = 23.58
$ws.Range("F587").Value = -36
$ws.Range("G587").Value = -710.28

# Rows 599,600: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B599").Value = 45709
$ws.Range("C599").Value = "CRE-Cremica Oatmeal Digestive 112.5 Gm"
$ws.Range("D599").Value = 13.15
$ws.Range("E599").Value = 15.69
$ws.Range("F599").Value = -300
$ws.Range("G599").Value = -3945
$ws.Range("B600").Value = 64925
$ws.Range("C600").Value = "CRE-Cremica Oatmeal Digestive 112.5 Gm"
$ws.Range("D600").Value = 13.15
$ws.Range("E600").Value = 13.97
$ws.Range("F600").Value = 302
$ws.Range("G600").Value = 3971.3

# Rows 687,688: cyclic shift of B,C,D,E,F,G (same product, different batches reordered)
$ws.Range("B687").Value = 64810
$ws.Range("C687").Value = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Range("D687").Value = 273.92
$ws.Range("E687").Value = 291.22
$ws.Range("F687").Value = 7
$ws.Range("G687").Value = 1917.44
$ws.Range("B688").Value = 53319
$ws.Range("C688").Value = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Range("D688").Value = 273.92
$ws.Range("E688").Value = 310.64
$ws.Range("F688").Value = -6
$ws.Range("G688").Value = -1643.52
